# Remove the top-ranked player (Nathan May, row 2) from the standings.
# This shifts every subsequent player up by one row and the sheet's
# used range shrinks from A1:J21 to A1:J20.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(2).Delete()

# After the shift, the "Rang" (rank) column (A) needs to be recomputed
# to reflect the new standings order (ties share the same rank).
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 11
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 19
